$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Already registered with the email suchitra@gmail.com, try to login"
$ws.Range("A3").Value = "Already registered with the email suchitra1@gmail.com, try to login"
$ws.Range("A4").Value = "Already registered with the email suchitra2@gmail.com, try to login"
$ws.Range("A5").Value = "Already registered with the email suchitra4@gmail.com, try to login"
